$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 53.50391
$ws.Range("H2").Value = 160.51173
$ws.Range("I2").Value = 0.2378988688891082
$ws.Range("J2").Value = 0.2378988688891082
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 22.93942688513
$ws.Range("R2").Value = 206.45484196617
$ws.Range("S2").Value = 0.0009802632286325806
$ws.Range("T2").Value = 0.0009802632286325806
$ws.Range("G3").Value = 53.50391
$ws.Range("H3").Value = 160.51173
$ws.Range("I3").Value = 0.2378988688891082
$ws.Range("J3").Value = 0.2378988688891082
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 4292.417364086669
$ws.Range("R3").Value = 38631.75627678003
$ws.Range("S3").Value = 0.1834265051619795
$ws.Range("T3").Value = 0.1834265051619795
$ws.Range("G4").Value = 53.50391
$ws.Range("H4").Value = 160.51173
$ws.Range("I4").Value = 0.2378988688891082
$ws.Range("J4").Value = 0.2378988688891082
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 1251.784309025843
$ws.Range("R4").Value = 11266.05878123259
$ws.Range("S4").Value = 0.05349210049849609
$ws.Range("T4").Value = 0.05349210049849609
$ws.Range("G5").Value = 82.26089466666666
$ws.Range("I5").Value = 0.3657634328968867
$ws.Range("J5").Value = 0.3657634328968868
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 35.26878276207066
$ws.Range("R5").Value = 317.419044858636
$ws.Range("S5").Value = 0.0015071296695167
$ws.Range("T5").Value = 0.001507129669516701
$ws.Range("G6").Value = 82.26089466666666
$ws.Range("I6").Value = 0.3657634328968867
$ws.Range("J6").Value = 0.3657634328968868
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 6599.482031360034
$ws.Range("R6").Value = 59395.33828224031
$ws.Range("S6").Value = 0.2820135653675476
$ws.Range("T6").Value = 0.2820135653675476
$ws.Range("G7").Value = 82.26089466666666
$ws.Range("I7").Value = 0.3657634328968867
$ws.Range("J7").Value = 0.3657634328968868
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 1924.586393595552
$ws.Range("R7").Value = 17321.27754235997
$ws.Range("S7").Value = 0.08224273785982246
$ws.Range("T7").Value = 0.08224273785982247
$ws.Range("G8").Value = 89.13710533333334
$ws.Range("H8").Value = 267.411316
$ws.Range("I8").Value = 0.3963376982140051
$ws.Range("J8").Value = 0.3963376982140052
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 38.21690995192934
$ws.Range("R8").Value = 343.952189567364
$ws.Range("S8").Value = 0.001633111050482399
$ws.Range("T8").Value = 0.0016331110504824
$ws.Range("G9").Value = 89.13710533333334
$ws.Range("H9").Value = 267.411316
$ws.Range("I9").Value = 0.3963376982140051
$ws.Range("J9").Value = 0.3963376982140052
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 7151.134537965963
$ws.Range("R9").Value = 64360.21084169367
$ws.Range("S9").Value = 0.3055871563694799
$ws.Range("T9").Value = 0.3055871563694799
$ws.Range("G10").Value = 89.13710533333334
$ws.Range("H10").Value = 267.411316
$ws.Range("I10").Value = 0.3963376982140051
$ws.Range("J10").Value = 0.3963376982140052
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 2085.463096215781
$ws.Range("R10").Value = 18769.16786594203
$ws.Range("S10").Value = 0.08911743079404287
$ws.Range("T10").Value = 0.08911743079404288
